# Update the "Comparison" sheet's column headers to the new virus labels
# (the underlying data columns / rows are unchanged - only the header text).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparison")

$ws.Range("B1").Value = "SARS-CoV2"
$ws.Range("C1").Value = "HRV"

# Move the active selection to E8 (matches the saved selection in the file).
$ws.Range("E8").Select()
